# Generate Report for Handoff
# The "b.md" file has been handed off again: update its status rows across
# the Overview, zh-cn and de-de sheets to reflect the new handoff, and
# record the error detail describing that the handback is stale.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/72c5f90b0544f077c5d6ed06f0bed88a042a6584/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cb73f45ec05a8dddaffc7256764df0fa9d92879c/e2e/b.md."

# --- Overview sheet : row for b.md (row 3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-01 20:43:47"

# --- zh-cn sheet : row for b.md (row 3) ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("F3").Value = "'False"
$wsZh.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-09-01 20:43:43"
$wsZh.Range("P3").Value = $errorDetail
$wsZh.Columns.Item(16).ColumnWidth = 40

# --- de-de sheet : row for b.md (row 3) ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("F3").Value = "'False"
$wsDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDe.Range("H3").Value = "2016-09-01 20:43:47"
$wsDe.Range("P3").Value = $errorDetail
$wsDe.Columns.Item(16).ColumnWidth = 40
